$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.961.25'
$ws.Range('E2').Value = '  -3.17%  '
$ws.Range('D3').Value = '2.289.81'
$ws.Range('E3').Value = '  -3.53%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.12%  '
$ws.Range('D5').Value = "'533.62"
$ws.Range('E5').Value = '  -4.12%  '
$ws.Range('D6').Value = "'131.01"
$ws.Range('E6').Value = '  -1.95%  '
$ws.Range('D7').Value = "'0.999"
$ws.Range('E7').Value = '  -0.10%  '
$ws.Range('D8').Value = "'0.583"
$ws.Range('E8').Value = '  -0.36%  '
$ws.Range('D9').Value = '2.287.87'
$ws.Range('E9').Value = '  -3.54%  '
$ws.Range('D10').Value = "'0.0995"
$ws.Range('E10').Value = '  -5.94%  '
$ws.Range('D11').Value = "'5.43"
$ws.Range('E11').Value = '  -3.97%  '
$ws.Range('E12').Value = '  -0.34%  '
$ws.Range('D13').Value = "'0.330"
$ws.Range('E13').Value = '  -3.68%  '
$ws.Range('D14').Value = "'23.50"
$ws.Range('E14').Value = '  -3.71%  '
$ws.Range('D15').Value = '2.697.81'
$ws.Range('E15').Value = '  -3.54%  '
$ws.Range('D16').Value = '57.901.07'
$ws.Range('E16').Value = '  -3.18%  '
$ws.Range('D17').Value = "'0.0000131"
$ws.Range('E17').Value = '  -4.77%  '
$ws.Range('D18').Value = '2.295.99'
$ws.Range('E18').Value = '  -3.30%  '
$ws.Range('D19').Value = "'10.49"
$ws.Range('E19').Value = '  -5.64%  '
$ws.Range('D20').Value = "'4.22"
$ws.Range('E20').Value = '  -5.87%  '
$ws.Range('D21').Value = "'311.95"
$ws.Range('E21').Value = '  -2.82%  '
$ws.Range('D22').Value = "'6.37"
$ws.Range('E22').Value = '  -4.11%  '
$ws.Range('E23').Value = '  +0.21%  '
$ws.Range('D24').Value = "'62.37"
$ws.Range('E24').Value = '  -2.71%  '
$ws.Range('E25').Value = '  -3.35%  '
$ws.Range('D26').Value = "'0.999"
$ws.Range('E26').Value = '  -0.16%  '
$ws.Range('D27').Value = "'7.99"
$ws.Range('E27').Value = '  -5.28%  '
$ws.Range('E28').Value = '  -6.84%  '
$ws.Range('D29').Value = "'170.67"
$ws.Range('E29').Value = '  +0.55%  '
$ws.Range('D30').Value = "'1.69"
$ws.Range('E30').Value = '  -6.21%  '
$ws.Range('D31').Value = '0.0₃0717'
$ws.Range('E31').Value = '  -5.53%  '
$ws.Range('D32').Value = "'5.74"
$ws.Range('E32').Value = '  -5.33%  '
$ws.Range('E33').Value = '  -7.05%  '
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D36').Value = "'17.71"
$ws.Range('E36').Value = '  -2.30%  '
$ws.Range('E37').Value = '  +0.02%  '
$ws.Range('E38').Value = '  -7.06%  '
$ws.Range('D39').Value = "'3.88"
$ws.Range('E39').Value = '  -6.10%  '
$ws.Range('D40').Value = "'38.14"
$ws.Range('E40').Value = '  -1.10%  '
$ws.Range('E41').Value = '  -6.47%  '
$ws.Range('D42').Value = "'141.60"
$ws.Range('E42').Value = '  -2.12%  '
$ws.Range('D43').Value = "'288.29"
$ws.Range('E43').Value = '  -9.38%  '
$ws.Range('D44').Value = "'3.41"
$ws.Range('E44').Value = '  -3.35%  '
$ws.Range('D45').Value = "'0.0945"
$ws.Range('E45').Value = '  -2.56%  '
$ws.Range('D46').Value = "'0.0495"
$ws.Range('E46').Value = '  -3.30%  '
$ws.Range('D47').Value = "'0.553"
$ws.Range('E47').Value = '  -2.55%  '
$ws.Range('D48').Value = "'18.08"
$ws.Range('E48').Value = '  -8.11%  '
$ws.Range('D49').Value = "'0.0210"
$ws.Range('E49').Value = '  -3.53%  '
$ws.Range('D50').Value = "'10.93"
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('E51').Value = '  -0.71%  '
